# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload"
#
# Adds a new day-16 entry (row 30) plus a block of follow-up notes
# (rows 31-43) describing a change request received from Mohan san on
# 25.01.2022, and moves the active selection down to the new notes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 30: new "No=16" entry --------------------------------------
$ws.Range("A30").Value = 16

$ws.Range("B30").Value = 44588
# "mm-dd-yy" is recognised as Excel's built-in short-date numFmtId (14),
# which is also how B2:B29 above are formatted (displays as m/d/yyyy).
$ws.Range("B30").NumberFormat = "mm-dd-yy"

$ws.Range("C30").Value = "RPA GSS"
$ws.Range("D30").Value = "We have received the change request for the customization parts at GSPN tasks as  follows on 25.01.2022 from Mohan san:"

# ---- Rows 31-43: bordered blank block (same look as rows above) -----
# Copy the plain-bordered style (s="2", used e.g. by A3) across the new
# block first, then fill in the specific cells - this reuses the
# existing border/style definitions instead of creating new ones.
$ws.Range("A3").Copy()
$ws.Range("A31:G43").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D31").Value = "1. We have to fetch the data and upload from the day was working till previous day  if the previous day is holiday."
$ws.Range("D32").Value = "2. if holidays are more than 3 days or more, then we have to split days and download and upload with date splitting "
$ws.Range("D33").Value = "3. if record counts is more than 5000, then again the divide the date or split the date and do the download and upload"

# Row 35: note + % complete + status, mirroring the pattern used by the
# other "No" rows above (percent cell keeps the percent+border style).
$ws.Range("D35").Value = "Note: point 1. has been  completed, whereas the testing has been done with  Warranty daily task, and other tasks are pending."

$ws.Range("E2").Copy()
$ws.Range("E35").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E35").Value = 0.2

$ws.Range("F35").Value = "WIP"

# ---- View state: selection moves to the new note at D31 -------------
$ws.Range("D31").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 19 | Out-Null
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
